$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A23").Value = "Longest Increasing Subsequence"
$ws.Range("B23").Value = "Return length of longest increasing subsequence"
$ws.Range("C23").Value = "Create a DP array filled with 1. Iterate over input array from left-right using 2 for loops. If n[i] > n[j] -> dp[i] = Math.max(dp[i], dp[j] + 1);   "
$ws.Range("D23").Value = "https://leetcode.com/problems/longest-increasing-subsequence/"

$ws.Hyperlinks.Add($ws.Range("D23"), "https://leetcode.com/problems/longest-increasing-subsequence/") | Out-Null

$ws.Range("A22:D22").Copy() | Out-Null
$ws.Range("A23:D23").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws.Range("C24").Select()
